$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A63 needs to hold the text "01-04-2021" without being auto-converted to a
# date serial by the COM Value-setter's literal-entry heuristics. Writing it
# as a formula that evaluates to the literal string, then collapsing the
# formula to its value via copy / paste-values, sidesteps that heuristic
# (and avoids minting a spurious new cell style in the process).
$r = $ws.Range("A63")
$r.Formula = '="01-04-2021"'
$r.Copy()
$r.PasteSpecial(-4163) # xlPasteValues

$ws.Cells.Item(63, 2).Value = 33.1
$ws.Cells.Item(63, 3).Value = 15.6
$ws.Cells.Item(63, 4).Value = 18.4
$ws.Cells.Item(63, 5).Value = -2.9
$ws.Cells.Item(63, 6).Value = 51.3
$ws.Cells.Item(63, 7).Value = 12.7
